$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 No -> Partially ; E2 (blank) -> "No Card draw message?"
$ws.Range("C2").Value = "Partially"
$ws.Range("E2").Value = "No Card draw message?"

# Row 3: B3 Partially -> Fully ; C3 No -> Partially ; E3 "Draw works" -> "No Card draw message?" (+ grey fill like E4)
$ws.Range("B3").Value = "Fully"
$ws.Range("C3").Value = "Partially"
$ws.Range("E4").Copy($ws.Range("E3"))
$ws.Range("E3").Value = "No Card draw message?"

# Row 4: B4 Partially -> Fully ; C4 No -> Partially ; E4 "Draw works" -> "No Card draw message?"
$ws.Range("B4").Value = "Fully"
$ws.Range("C4").Value = "Partially"
$ws.Range("E4").Value = "No Card draw message?"

# Row 5: B5 No -> Fully ; C5 No -> Partially ; E5 (new) -> "No Card draw message?" (grey fill like E4)
$ws.Range("B5").Value = "Fully"
$ws.Range("C5").Value = "Partially"
$ws.Range("E4").Copy($ws.Range("E5"))
$ws.Range("E5").Value = "No Card draw message?"

# Update the active selection to match the saved view state
$ws.Range("D33").Select()
